$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Every Price/Volume cell is stored as text in the source sheet (t="inlineStr").
# Excel's COM Value setter auto-infers numeric-looking strings as numbers, so we
# force text interpretation with NumberFormat "@" and then restore the default
# "Normal" style so no stray style index is left on the cell.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "41.616.04"
Set-TextValue $ws.Range("E2") "  +0.33%  "

Set-TextValue $ws.Range("D3") "2.470.63"
Set-TextValue $ws.Range("E3") "  +0.15%  "

Set-TextValue $ws.Range("D4") "1.00"
Set-TextValue $ws.Range("E4") "  -0.86%  "

Set-TextValue $ws.Range("D5") "314.41"
Set-TextValue $ws.Range("E5") "  +0.68%  "

Set-TextValue $ws.Range("D6") "91.45"
Set-TextValue $ws.Range("E6") "  +0.39%  "

Set-TextValue $ws.Range("D7") "0.550"
Set-TextValue $ws.Range("E7") "  +1.81%  "

Set-TextValue $ws.Range("E8") "  -0.89%  "

Set-TextValue $ws.Range("D9") "0.513"
Set-TextValue $ws.Range("E9") "  +4.81%  "

Set-TextValue $ws.Range("D10") "32.67"
Set-TextValue $ws.Range("E10") "  +0.16%  "

Set-TextValue $ws.Range("D11") "0.0794"
Set-TextValue $ws.Range("E11") "  +2.47%  "

Set-TextValue $ws.Range("E12") "  +0.69%  "

Set-TextValue $ws.Range("D13") "2.845.35"
Set-TextValue $ws.Range("E13") "  -0.10%  "

Set-TextValue $ws.Range("D14") "6.88"
Set-TextValue $ws.Range("E14") "  +1.02%  "

Set-TextValue $ws.Range("D15") "15.91"
Set-TextValue $ws.Range("E15") "  +4.87%  "

Set-TextValue $ws.Range("D16") "2.452.54"
Set-TextValue $ws.Range("E16") "  -3.56%  "

Set-TextValue $ws.Range("D17") "0.778"
Set-TextValue $ws.Range("E17") "  +0.10%  "

Set-TextValue $ws.Range("D18") "41.585.47"
Set-TextValue $ws.Range("E18") "  +0.70%  "

Set-TextValue $ws.Range("D19") "6.51"
Set-TextValue $ws.Range("E19") "  +4.12%  "

Set-TextValue $ws.Range("E20") "  +3.12%  "

Set-TextValue $ws.Range("D21") "71.31"
Set-TextValue $ws.Range("E21") "  +1.50%  "

Set-TextValue $ws.Range("D22") "11.20"
Set-TextValue $ws.Range("E22") "  +2.62%  "

Set-TextValue $ws.Range("D23") "238.44"
Set-TextValue $ws.Range("E23") "  +1.87%  "

Set-TextValue $ws.Range("E24") "  +0.20%  "

Set-TextValue $ws.Range("E25") "  +2.85%  "

Set-TextValue $ws.Range("E26") "  -0.09%  "

Set-TextValue $ws.Range("D27") "24.67"
Set-TextValue $ws.Range("E27") "  +3.55%  "

Set-TextValue $ws.Range("D28") "2.26"
Set-TextValue $ws.Range("E28") "  +0.68%  "

Set-TextValue $ws.Range("D29") "9.67"
Set-TextValue $ws.Range("E29") "  +0.19%  "

Set-TextValue $ws.Range("D30") "35.40"
Set-TextValue $ws.Range("E30") "  -1.43%  "

Set-TextValue $ws.Range("D31") "156.37"
Set-TextValue $ws.Range("E31") "  +2.61%  "

Set-TextValue $ws.Range("D32") "5.45"
Set-TextValue $ws.Range("E32") "  +1.25%  "

Set-TextValue $ws.Range("E33") "  +1.02%  "

Set-TextValue $ws.Range("E34") "  +1.85%  "

Set-TextValue $ws.Range("D35") "17.28"
Set-TextValue $ws.Range("E35") "  -0.38%  "

Set-TextValue $ws.Range("D36") "2.35"
Set-TextValue $ws.Range("E36") "  -8.62%  "

Set-TextValue $ws.Range("D37") "2.89"
Set-TextValue $ws.Range("E37") "  -2.65%  "

Set-TextValue $ws.Range("D38") "0.115"
Set-TextValue $ws.Range("E38") "  +1.92%  "

Set-TextValue $ws.Range("D39") "0.103"
Set-TextValue $ws.Range("E39") "  +4.15%  "

Set-TextValue $ws.Range("E40") "  -1.83%  "

Set-TextValue $ws.Range("D41") "4.00"
Set-TextValue $ws.Range("E41") "  -0.38%  "

Set-TextValue $ws.Range("E42") "  -1.20%  "

Set-TextValue $ws.Range("D43") "1.960.99"
Set-TextValue $ws.Range("E43") "  +0.28%  "

Set-TextValue $ws.Range("E44") "  +1.13%  "

Set-TextValue $ws.Range("D45") "18.55"
Set-TextValue $ws.Range("E45") "  -3.11%  "

Set-TextValue $ws.Range("D46") "2.92"
Set-TextValue $ws.Range("E46") "  -0.10%  "

Set-TextValue $ws.Range("D47") "9.05"
Set-TextValue $ws.Range("E47") "  +5.09%  "

Set-TextValue $ws.Range("D48") "2.703.43"
Set-TextValue $ws.Range("E48") "  -0.45%  "

Set-TextValue $ws.Range("D49") "97.21"
Set-TextValue $ws.Range("E49") "  +2.27%  "

Set-TextValue $ws.Range("D50") "67.33"
Set-TextValue $ws.Range("E50") "  -0.05%  "

Set-TextValue $ws.Range("B51") "Algorand"
Set-TextValue $ws.Range("C51") "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws.Range("D51") "0.172"
Set-TextValue $ws.Range("E51") "  -1.15%  "
